$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# D-column values are textual (e.g. thousands separated with dots, fixed decimals)
# so force text format to avoid Excel auto-converting them to numbers and losing
# formatting (trailing zeros, multiple-dot thousands grouping), then restore the
# default "Normal" style so no stray per-cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.473.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.591.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.06%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.588.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.257.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.367.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "

$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.578.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.489"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "510.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.90"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.780.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.98%  "

$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.146"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.565"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "602.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.24%  "

$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.918"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "35.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("E51").Value = "  -0.73%  "
